$wb = $excel.ActiveWorkbook
$wsTypo = $wb.Worksheets.Item("Typography")
$wsTrans = $wb.Worksheets.Item("Translation")

# --- Typography sheet: append the degree sign to the "Wildcard Characters" column (G) ---
$wsTypo.Cells.Item(4, 7).Value = '.",°'
$wsTypo.Cells.Item(5, 7).Value = '.",°'
$wsTypo.Cells.Item(6, 7).Value = '.",_°'
$wsTypo.Cells.Item(7, 7).Value = '.",°'
$wsTypo.Cells.Item(8, 7).Value = '.",°'

# --- Translation sheet: add new rows 57-90 (GRBL laser control texts) ---
$wsTrans.Cells.Item(57, 2).Value = 'SingleUseId63'
$wsTrans.Cells.Item(57, 3).Value = 'Small'
$wsTrans.Cells.Item(57, 4).Value = 'Left'
$wsTrans.Cells.Item(57, 5).Value = 'LTR'
$wsTrans.Cells.Item(57, 6).Value = ' <uartConsoleBfr>'
$wsTrans.Cells.Item(57, 7).Value = ' <uartConsoleBfr_0>'
$wsTrans.Cells.Item(58, 2).Value = 'SingleUseId64'
$wsTrans.Cells.Item(58, 3).Value = 'Small'
$wsTrans.Cells.Item(58, 4).Value = 'Left'
$wsTrans.Cells.Item(58, 5).Value = 'LTR'
$wsTrans.Cells.Item(58, 6).Value = ' <uartConsoleBfr>'
$wsTrans.Cells.Item(58, 7).Value = ' <uartConsoleBfr_1>'
$wsTrans.Cells.Item(59, 2).Value = 'SingleUseId65'
$wsTrans.Cells.Item(59, 3).Value = 'Small'
$wsTrans.Cells.Item(59, 4).Value = 'Left'
$wsTrans.Cells.Item(59, 5).Value = 'LTR'
$wsTrans.Cells.Item(59, 6).Value = ' <uartConsoleBfr>'
$wsTrans.Cells.Item(59, 7).Value = ' <uartConsoleBfr_2>'
$wsTrans.Cells.Item(60, 2).Value = 'SingleUseId66'
$wsTrans.Cells.Item(60, 3).Value = 'Small'
$wsTrans.Cells.Item(60, 4).Value = 'Left'
$wsTrans.Cells.Item(60, 5).Value = 'LTR'
$wsTrans.Cells.Item(60, 6).Value = ' <uartConsoleBfr>'
$wsTrans.Cells.Item(60, 7).Value = ' <uartConsoleBfr_3>'
$wsTrans.Cells.Item(61, 2).Value = 'SingleUseId67'
$wsTrans.Cells.Item(61, 3).Value = 'Small'
$wsTrans.Cells.Item(61, 4).Value = 'Left'
$wsTrans.Cells.Item(61, 5).Value = 'LTR'
$wsTrans.Cells.Item(61, 6).Value = ' <uartConsoleBfr>'
$wsTrans.Cells.Item(61, 7).Value = ' <uartConsoleBfr_4>'
$wsTrans.Cells.Item(62, 2).Value = 'SingleUseId68'
$wsTrans.Cells.Item(62, 3).Value = 'Small'
$wsTrans.Cells.Item(62, 4).Value = 'Left'
$wsTrans.Cells.Item(62, 5).Value = 'LTR'
$wsTrans.Cells.Item(62, 6).Value = ' <uartConsoleBfr>'
$wsTrans.Cells.Item(62, 7).Value = ' <uartConsoleBfr_5>'
$wsTrans.Cells.Item(63, 2).Value = 'SingleUseId69'
$wsTrans.Cells.Item(63, 3).Value = 'Large'
$wsTrans.Cells.Item(63, 4).Value = 'Right'
$wsTrans.Cells.Item(63, 5).Value = 'LTR'
$wsTrans.Cells.Item(63, 6).Value = 'alpha: °'
$wsTrans.Cells.Item(63, 7).Value = 'alpha [°]:'
$wsTrans.Cells.Item(64, 2).Value = 'SingleUseId71'
$wsTrans.Cells.Item(64, 3).Value = 'Large'
$wsTrans.Cells.Item(64, 4).Value = 'Right'
$wsTrans.Cells.Item(64, 5).Value = 'LTR'
$wsTrans.Cells.Item(64, 6).Value = 'širina [mm]:'
$wsTrans.Cells.Item(64, 7).Value = 'width [mm]:'
$wsTrans.Cells.Item(65, 2).Value = 'SingleUseId72'
$wsTrans.Cells.Item(65, 3).Value = 'Default'
$wsTrans.Cells.Item(65, 4).Value = 'Right'
$wsTrans.Cells.Item(65, 5).Value = 'LTR'
$wsTrans.Cells.Item(65, 6).Value = 'hitrost [mm/s]:'
$wsTrans.Cells.Item(65, 7).Value = 'f. rate [mm/s]:'
$wsTrans.Cells.Item(66, 2).Value = 'SingleUseId73'
$wsTrans.Cells.Item(66, 3).Value = 'Large'
$wsTrans.Cells.Item(66, 4).Value = 'Right'
$wsTrans.Cells.Item(66, 5).Value = 'LTR'
$wsTrans.Cells.Item(66, 6).Value = 'beta [°]:'
$wsTrans.Cells.Item(66, 7).Value = 'beta [°]:'
$wsTrans.Cells.Item(67, 2).Value = 'SingleUseId74'
$wsTrans.Cells.Item(67, 3).Value = 'Large'
$wsTrans.Cells.Item(67, 4).Value = 'Left'
$wsTrans.Cells.Item(67, 5).Value = 'LTR'
$wsTrans.Cells.Item(67, 6).Value = '.'
$wsTrans.Cells.Item(67, 7).Value = '.'
$wsTrans.Cells.Item(68, 2).Value = 'SingleUseId75'
$wsTrans.Cells.Item(68, 3).Value = 'Large'
$wsTrans.Cells.Item(68, 4).Value = 'Left'
$wsTrans.Cells.Item(68, 5).Value = 'LTR'
$wsTrans.Cells.Item(68, 6).Value = '.'
$wsTrans.Cells.Item(68, 7).Value = '.'
$wsTrans.Cells.Item(69, 2).Value = 'SingleUseId76'
$wsTrans.Cells.Item(69, 3).Value = 'Large'
$wsTrans.Cells.Item(69, 4).Value = 'Left'
$wsTrans.Cells.Item(69, 5).Value = 'LTR'
$wsTrans.Cells.Item(69, 6).Value = '.'
$wsTrans.Cells.Item(69, 7).Value = '.'
$wsTrans.Cells.Item(70, 2).Value = 'SingleUseId77'
$wsTrans.Cells.Item(70, 3).Value = 'Large'
$wsTrans.Cells.Item(70, 4).Value = 'Left'
$wsTrans.Cells.Item(70, 5).Value = 'LTR'
$wsTrans.Cells.Item(70, 6).Value = '.'
$wsTrans.Cells.Item(70, 7).Value = '.'
$wsTrans.Cells.Item(71, 2).Value = 'SingleUseId78'
$wsTrans.Cells.Item(71, 3).Value = 'Large'
$wsTrans.Cells.Item(71, 4).Value = 'Left'
$wsTrans.Cells.Item(71, 5).Value = 'LTR'
$wsTrans.Cells.Item(71, 6).Value = '.'
$wsTrans.Cells.Item(71, 7).Value = '.'
$wsTrans.Cells.Item(72, 2).Value = 'SingleUseId79'
$wsTrans.Cells.Item(72, 3).Value = 'Large'
$wsTrans.Cells.Item(72, 4).Value = 'Left'
$wsTrans.Cells.Item(72, 5).Value = 'LTR'
$wsTrans.Cells.Item(72, 6).Value = '.'
$wsTrans.Cells.Item(72, 7).Value = '.'
$wsTrans.Cells.Item(73, 2).Value = 'SingleUseId80'
$wsTrans.Cells.Item(73, 3).Value = 'Default'
$wsTrans.Cells.Item(73, 4).Value = 'Right'
$wsTrans.Cells.Item(73, 5).Value = 'LTR'
$wsTrans.Cells.Item(73, 6).Value = 'alpha: [°]'
$wsTrans.Cells.Item(73, 7).Value = 'alpha [°]:'
$wsTrans.Cells.Item(74, 2).Value = 'SingleUseId81'
$wsTrans.Cells.Item(74, 3).Value = 'Default'
$wsTrans.Cells.Item(74, 4).Value = 'Right'
$wsTrans.Cells.Item(74, 5).Value = 'LTR'
$wsTrans.Cells.Item(74, 6).Value = 'širina [mm]:'
$wsTrans.Cells.Item(74, 7).Value = 'width [mm]:'
$wsTrans.Cells.Item(75, 2).Value = 'SingleUseId82'
$wsTrans.Cells.Item(75, 3).Value = 'Default'
$wsTrans.Cells.Item(75, 4).Value = 'Right'
$wsTrans.Cells.Item(75, 5).Value = 'LTR'
$wsTrans.Cells.Item(75, 6).Value = 'hitrost [mm/s]:'
$wsTrans.Cells.Item(75, 7).Value = 'f. rate [mm/s]:'
$wsTrans.Cells.Item(76, 2).Value = 'SingleUseId83'
$wsTrans.Cells.Item(76, 3).Value = 'Default'
$wsTrans.Cells.Item(76, 4).Value = 'Right'
$wsTrans.Cells.Item(76, 5).Value = 'LTR'
$wsTrans.Cells.Item(76, 6).Value = 'beta [°]:'
$wsTrans.Cells.Item(76, 7).Value = 'beta [°]:'
$wsTrans.Cells.Item(77, 2).Value = 'SingleUseId84'
$wsTrans.Cells.Item(77, 3).Value = 'Default'
$wsTrans.Cells.Item(77, 4).Value = 'Center'
$wsTrans.Cells.Item(77, 5).Value = 'LTR'
$wsTrans.Cells.Item(77, 6).Value = 'Save 1'
$wsTrans.Cells.Item(77, 7).Value = 'Save 1'
$wsTrans.Cells.Item(78, 2).Value = 'SingleUseId85'
$wsTrans.Cells.Item(78, 3).Value = 'Default'
$wsTrans.Cells.Item(78, 4).Value = 'Center'
$wsTrans.Cells.Item(78, 5).Value = 'LTR'
$wsTrans.Cells.Item(78, 6).Value = 'Save 2'
$wsTrans.Cells.Item(78, 7).Value = 'Save 2'
$wsTrans.Cells.Item(79, 2).Value = 'SingleUseId86'
$wsTrans.Cells.Item(79, 3).Value = 'Default'
$wsTrans.Cells.Item(79, 4).Value = 'Center'
$wsTrans.Cells.Item(79, 5).Value = 'LTR'
$wsTrans.Cells.Item(79, 6).Value = 'Load 1'
$wsTrans.Cells.Item(79, 7).Value = 'Load 1'
$wsTrans.Cells.Item(80, 2).Value = 'SingleUseId87'
$wsTrans.Cells.Item(80, 3).Value = 'Default'
$wsTrans.Cells.Item(80, 4).Value = 'Center'
$wsTrans.Cells.Item(80, 5).Value = 'LTR'
$wsTrans.Cells.Item(80, 6).Value = 'Load 2'
$wsTrans.Cells.Item(80, 7).Value = 'Load 2'
$wsTrans.Cells.Item(81, 2).Value = 'SingleUseId88'
$wsTrans.Cells.Item(81, 3).Value = 'Default'
$wsTrans.Cells.Item(81, 4).Value = 'Center'
$wsTrans.Cells.Item(81, 5).Value = 'LTR'
$wsTrans.Cells.Item(81, 6).Value = 'Save 3'
$wsTrans.Cells.Item(81, 7).Value = 'Save 3'
$wsTrans.Cells.Item(82, 2).Value = 'SingleUseId89'
$wsTrans.Cells.Item(82, 3).Value = 'Default'
$wsTrans.Cells.Item(82, 4).Value = 'Center'
$wsTrans.Cells.Item(82, 5).Value = 'LTR'
$wsTrans.Cells.Item(82, 6).Value = 'Load 3'
$wsTrans.Cells.Item(82, 7).Value = 'Load 3'
$wsTrans.Cells.Item(83, 2).Value = 'SingleUseId90'
$wsTrans.Cells.Item(83, 3).Value = 'Large'
$wsTrans.Cells.Item(83, 4).Value = 'Left'
$wsTrans.Cells.Item(83, 5).Value = 'LTR'
$wsTrans.Cells.Item(83, 6).Value = 'alpha: <alpha>°'
$wsTrans.Cells.Item(83, 7).Value = 'alpha: <alpha>°'
$wsTrans.Cells.Item(84, 2).Value = 'SingleUseId91'
$wsTrans.Cells.Item(84, 3).Value = 'Large'
$wsTrans.Cells.Item(84, 4).Value = 'Left'
$wsTrans.Cells.Item(84, 5).Value = 'LTR'
$wsTrans.Cells.Item(84, 6).Value = 'širina: <width>mm'
$wsTrans.Cells.Item(84, 7).Value = 'width: <width>mm'
$wsTrans.Cells.Item(85, 2).Value = 'SingleUseId92'
$wsTrans.Cells.Item(85, 3).Value = 'Large'
$wsTrans.Cells.Item(85, 4).Value = 'Left'
$wsTrans.Cells.Item(85, 5).Value = 'LTR'
$wsTrans.Cells.Item(85, 6).Value = 'hitrost: <feedrate>mm/s'
$wsTrans.Cells.Item(85, 7).Value = 'f. rate: <feedrate>mm/s'
$wsTrans.Cells.Item(86, 2).Value = 'SingleUseId93'
$wsTrans.Cells.Item(86, 3).Value = 'Large'
$wsTrans.Cells.Item(86, 4).Value = 'Left'
$wsTrans.Cells.Item(86, 5).Value = 'LTR'
$wsTrans.Cells.Item(86, 6).Value = 'beta: <beta>°'
$wsTrans.Cells.Item(86, 7).Value = 'beta: <beta>°'
$wsTrans.Cells.Item(87, 2).Value = 'SingleUseId94'
$wsTrans.Cells.Item(87, 3).Value = 'Large'
$wsTrans.Cells.Item(87, 4).Value = 'Left'
$wsTrans.Cells.Item(87, 5).Value = 'LTR'
$wsTrans.Cells.Item(87, 6).Value = '.'
$wsTrans.Cells.Item(87, 7).Value = '.'
$wsTrans.Cells.Item(88, 2).Value = 'SingleUseId96'
$wsTrans.Cells.Item(88, 3).Value = 'Large'
$wsTrans.Cells.Item(88, 4).Value = 'Left'
$wsTrans.Cells.Item(88, 5).Value = 'LTR'
$wsTrans.Cells.Item(88, 6).Value = '.'
$wsTrans.Cells.Item(88, 7).Value = '.'
$wsTrans.Cells.Item(89, 2).Value = 'SingleUseId98'
$wsTrans.Cells.Item(89, 3).Value = 'Large'
$wsTrans.Cells.Item(89, 4).Value = 'Right'
$wsTrans.Cells.Item(89, 5).Value = 'LTR'
$wsTrans.Cells.Item(89, 6).Value = 'x0 [mm]:'
$wsTrans.Cells.Item(89, 7).Value = 'x0 [mm]:'
$wsTrans.Cells.Item(90, 2).Value = 'SingleUseId100'
$wsTrans.Cells.Item(90, 3).Value = 'Large'
$wsTrans.Cells.Item(90, 4).Value = 'Right'
$wsTrans.Cells.Item(90, 5).Value = 'LTR'
$wsTrans.Cells.Item(90, 6).Value = 'y0 [mm]:'
$wsTrans.Cells.Item(90, 7).Value = 'y0 [mm]:'
